$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 2999289
$ws.Range("B2").Value = 96319
$ws.Range("C2").Value = "Ovaliderad"
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 219799
$ws.Range("F2").Value = "Kärrknipprot"
$ws.Range("G2").Value = "Epipactis palustris"
$ws.Range("H2").Value = "(L.) Crantz"
$ws.Range("M2").ClearContents()
$ws.Range("P2").Value = "06J1A05, Gtl"
$ws.Range("Q2").Value = 698621.4564279296
$ws.Range("R2").Value = 6356902.266207782
$ws.Range("S2").Value = 50
$ws.Range("T2").Value = "Gotland"
$ws.Range("U2").Value = "Gotland"
$ws.Range("V2").Value = "Gotland"
$ws.Range("W2").Value = "Gerum"
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2007-07-09"
$ws.Range("Z2").NumberFormat = "@"
$ws.Range("Z2").Value = "00:00"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2007-09-05"
$ws.Range("AB2").NumberFormat = "@"
$ws.Range("AB2").Value = "00:00"
$ws.Range("AC2").Value = "Lokalens storlek varierar starkt från 5 m i radie till flera hundra meter. Valde 50 m i denna rapportering.Mittkoordinater för varje rikkärrsobjekt, ej för varje observation."
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AG2").Value = $false
$ws.Range("AH2").Value = "Rikkärr"
$ws.Range("AW2").Value = "Britta Johansson"
$ws.Range("AX2").Value = "Magnus Martinsson"
$ws.Range("AY2").Value = "Åtgärdsprogram för hotade arter"

# Row 3
$ws.Range("A3").Value = 5168355
$ws.Range("B3").Value = 97335
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 222662
$ws.Range("F3").Value = "Axag"
$ws.Range("G3").Value = "Schoenus ferrugineus"
$ws.Range("H3").Value = "L."
$ws.Range("M3").ClearContents()
$ws.Range("P3").Value = "06J1A05, Gtl"
$ws.Range("Q3").Value = 698621.4564279296
$ws.Range("R3").Value = 6356902.266207782
$ws.Range("S3").Value = 50
$ws.Range("T3").Value = "Gotland"
$ws.Range("U3").Value = "Gotland"
$ws.Range("V3").Value = "Gotland"
$ws.Range("W3").Value = "Gerum"
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2007-07-09"
$ws.Range("Z3").NumberFormat = "@"
$ws.Range("Z3").Value = "00:00"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2007-09-05"
$ws.Range("AB3").NumberFormat = "@"
$ws.Range("AB3").Value = "00:00"
$ws.Range("AC3").Value = "Lokalens storlek varierar starkt från 5 m i radie till flera hundra meter. Valde 50 m i denna rapportering.Mittkoordinater för varje rikkärrsobjekt, ej för varje observation."
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AH3").Value = "Rikkärr"
$ws.Range("AW3").Value = "Britta Johansson"
$ws.Range("AX3").Value = "Magnus Martinsson"
$ws.Range("AY3").Value = "Åtgärdsprogram för hotade arter"

# Row 4
$ws.Range("A4").Value = 2289844
$ws.Range("B4").Value = 96336
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 219811
$ws.Range("F4").Value = "Brudsporre"
$ws.Range("G4").Value = "Gymnadenia conopsea"
$ws.Range("H4").Value = "(L.) R. Br."
$ws.Range("M4").ClearContents()
$ws.Range("P4").Value = "06J1A05, Gtl"
$ws.Range("Q4").Value = 698621.4564279296
$ws.Range("R4").Value = 6356902.266207782
$ws.Range("S4").Value = 50
$ws.Range("T4").Value = "Gotland"
$ws.Range("U4").Value = "Gotland"
$ws.Range("V4").Value = "Gotland"
$ws.Range("W4").Value = "Gerum"
$ws.Range("Y4").NumberFormat = "@"
$ws.Range("Y4").Value = "2007-07-09"
$ws.Range("Z4").NumberFormat = "@"
$ws.Range("Z4").Value = "00:00"
$ws.Range("AA4").NumberFormat = "@"
$ws.Range("AA4").Value = "2007-09-05"
$ws.Range("AB4").NumberFormat = "@"
$ws.Range("AB4").Value = "00:00"
$ws.Range("AC4").Value = "Lokalens storlek varierar starkt från 5 m i radie till flera hundra meter. Valde 50 m i denna rapportering.Mittkoordinater för varje rikkärrsobjekt, ej för varje observation."
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AH4").Value = "Rikkärr"
$ws.Range("AW4").Value = "Britta Johansson"
$ws.Range("AX4").Value = "Magnus Martinsson"
$ws.Range("AY4").Value = "Åtgärdsprogram för hotade arter"

# Row 5
$ws.Range("A5").Value = 4095111
$ws.Range("B5").Value = 95990
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 221930
$ws.Range("F5").Value = "Kärrlilja"
$ws.Range("G5").Value = "Tofieldia calyculata"
$ws.Range("H5").Value = "(L.) Wahlenb."
$ws.Range("M5").ClearContents()
$ws.Range("P5").Value = "06J1A05, Gtl"
$ws.Range("Q5").Value = 698621.4564279296
$ws.Range("R5").Value = 6356902.266207782
$ws.Range("S5").Value = 50
$ws.Range("T5").Value = "Gotland"
$ws.Range("U5").Value = "Gotland"
$ws.Range("V5").Value = "Gotland"
$ws.Range("W5").Value = "Gerum"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2007-07-09"
$ws.Range("Z5").NumberFormat = "@"
$ws.Range("Z5").Value = "00:00"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2007-09-05"
$ws.Range("AB5").NumberFormat = "@"
$ws.Range("AB5").Value = "00:00"
$ws.Range("AC5").Value = "Lokalens storlek varierar starkt från 5 m i radie till flera hundra meter. Valde 50 m i denna rapportering.Mittkoordinater för varje rikkärrsobjekt, ej för varje observation."
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AH5").Value = "Rikkärr"
$ws.Range("AW5").Value = "Britta Johansson"
$ws.Range("AX5").Value = "Magnus Martinsson"
$ws.Range("AY5").Value = "Åtgärdsprogram för hotade arter"

# Row 6
$ws.Range("A6").Value = 3704166
$ws.Range("B6").Value = 103164
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 221137
$ws.Range("F6").Value = "Majviva"
$ws.Range("G6").Value = "Primula farinosa"
$ws.Range("H6").Value = "L."
$ws.Range("M6").ClearContents()
$ws.Range("P6").Value = "06J1A05, Gtl"
$ws.Range("Q6").Value = 698621.4564279296
$ws.Range("R6").Value = 6356902.266207782
$ws.Range("S6").Value = 50
$ws.Range("T6").Value = "Gotland"
$ws.Range("U6").Value = "Gotland"
$ws.Range("V6").Value = "Gotland"
$ws.Range("W6").Value = "Gerum"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2007-07-09"
$ws.Range("Z6").NumberFormat = "@"
$ws.Range("Z6").Value = "00:00"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2007-09-05"
$ws.Range("AB6").NumberFormat = "@"
$ws.Range("AB6").Value = "00:00"
$ws.Range("AC6").Value = "Lokalens storlek varierar starkt från 5 m i radie till flera hundra meter. Valde 50 m i denna rapportering.Mittkoordinater för varje rikkärrsobjekt, ej för varje observation."
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AH6").Value = "Rikkärr"
$ws.Range("AW6").Value = "Britta Johansson"
$ws.Range("AX6").Value = "Magnus Martinsson"
$ws.Range("AY6").Value = "Åtgärdsprogram för hotade arter"

# Row 7
$ws.Range("A7").Value = 94011938
$ws.Range("B7").Value = 98520
$ws.Range("C7").Value = "Ovaliderad"
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 222498
$ws.Range("F7").Value = "Blåsippa"
$ws.Range("G7").Value = "Hepatica nobilis"
$ws.Range("H7").Value = "Schreb."
$ws.Range("M7").ClearContents()
$ws.Range("P7").Value = "Ygne-Hemse, Gtl"
$ws.Range("Q7").Value = 698779.9316287825
$ws.Range("R7").Value = 6356965.133827931
$ws.Range("S7").Value = 10
$ws.Range("T7").Value = "Gotland"
$ws.Range("U7").Value = "Gotland"
$ws.Range("V7").Value = "Gotland"
$ws.Range("W7").Value = "Gerum"
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2021-05-24"
$ws.Range("Z7").NumberFormat = "@"
$ws.Range("Z7").Value = "00:00"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2021-05-24"
$ws.Range("AB7").NumberFormat = "@"
$ws.Range("AB7").Value = "00:00"
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false
$ws.Range("AH7").ClearContents()
$ws.Range("AW7").Value = "Tony Svensson"
$ws.Range("AX7").Value = "Tony Svensson"
$ws.Range("AY7").Value = "Ecogain"

# Row 8
$ws.Range("A8").Value = 94011937
$ws.Range("B8").Value = 98520
$ws.Range("C8").Value = "Ovaliderad"
$ws.Range("D8").Value = "LC"
$ws.Range("E8").Value = 222498
$ws.Range("F8").Value = "Blåsippa"
$ws.Range("G8").Value = "Hepatica nobilis"
$ws.Range("H8").Value = "Schreb."
$ws.Range("M8").ClearContents()
$ws.Range("P8").Value = "Ygne-Hemse, Gtl"
$ws.Range("Q8").Value = 698786.3346319427
$ws.Range("R8").Value = 6356944.887886292
$ws.Range("S8").Value = 10
$ws.Range("T8").Value = "Gotland"
$ws.Range("U8").Value = "Gotland"
$ws.Range("V8").Value = "Gotland"
$ws.Range("W8").Value = "Gerum"
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = "2021-05-24"
$ws.Range("Z8").NumberFormat = "@"
$ws.Range("Z8").Value = "00:00"
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = "2021-05-24"
$ws.Range("AB8").NumberFormat = "@"
$ws.Range("AB8").Value = "00:00"
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AG8").Value = $false
$ws.Range("AH8").ClearContents()
$ws.Range("AW8").Value = "Tony Svensson"
$ws.Range("AX8").Value = "Tony Svensson"
$ws.Range("AY8").Value = "Ecogain"

# Row 9
$ws.Range("A9").Value = 94011889
$ws.Range("B9").Value = 56887
$ws.Range("C9").Value = "Ovaliderad"
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 102995
$ws.Range("F9").Value = "Buskskvätta"
$ws.Range("G9").Value = "Saxicola rubetra"
$ws.Range("H9").Value = "(Linnaeus, 1758)"
$ws.Range("M9").Value = "spel/sång"
$ws.Range("P9").Value = "Ygne-Hemse, Gtl"
$ws.Range("Q9").Value = 698742.6195520113
$ws.Range("R9").Value = 6356963.322957435
$ws.Range("S9").Value = 10
$ws.Range("T9").Value = "Gotland"
$ws.Range("U9").Value = "Gotland"
$ws.Range("V9").Value = "Gotland"
$ws.Range("W9").Value = "Gerum"
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2021-05-24"
$ws.Range("Z9").NumberFormat = "@"
$ws.Range("Z9").Value = "00:00"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = "2021-05-24"
$ws.Range("AB9").NumberFormat = "@"
$ws.Range("AB9").Value = "00:00"
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").Value = $false
$ws.Range("AE9").Value = $false
$ws.Range("AG9").Value = $false
$ws.Range("AH9").ClearContents()
$ws.Range("AW9").Value = "Tony Svensson"
$ws.Range("AX9").Value = "Tony Svensson"
$ws.Range("AY9").Value = "Ecogain"
